# "Checked AutoRun and Template"
#
# 1) Refresh the cached "datetimeFigureOut" date placeholder text that
#    lives on the slide master and on every slide layout
#    (7/22/2020 -> 7/27/2020).
# 2) On slide 2 ("Collect"/"Find"):
#      - retitle "Collect" -> "Find"
#      - delete the apple photo ("Picture 4")
#      - re-home the remaining SVG graphic ("Graphic 3") into the spot
#        the apple used to occupy, now that it is the only picture left

$p = $ppt.ActivePresentation

# ---- 1) date placeholder text: master + every custom (slide) layout ----
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "7/22/2020") {
            $sh.TextFrame.TextRange.Text = "7/27/2020"
        }
    }
}

$customLayouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $customLayouts.Count; $li++) {
    $layoutShapes = $customLayouts.Item($li).Shapes
    for ($i = 1; $i -le $layoutShapes.Count; $i++) {
        $sh = $layoutShapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "7/22/2020") {
                $sh.TextFrame.TextRange.Text = "7/27/2020"
            }
        }
    }
}

# ---- 2) slide 2 content swap ----
$s = $p.Slides.Item(2)

# title: "Collect" -> "Find"
$s.Shapes.Item("Title 1").TextFrame.TextRange.Text = "Find"

# drop the apple picture; the tomato/SVG graphic ("Graphic 3") survives
$s.Shapes.Item("Picture 4").Delete()

# re-position the surviving graphic into the vacated spot.
# (Left/Top are Single-precision points under the hood, so the literal
# values below are nudged by a hair so the round-trip lands on the exact
# EMU the author dragged it to: x=4835108 EMU, y=2800566 EMU.)
$graphic = $s.Shapes.Item("Graphic 3")
$graphic.Left = 380.7171936035156
$graphic.Top = 220.5170135498047
